$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff": the report table re-sorts so the
# fbebf677... source file is listed before the 36038c19... source file, and
# the 36038c19... file's handoff status/timestamps/error detail get
# refreshed to reflect a brand-new handoff ("Ready for handoff").
# Only the cells whose text actually changes are touched below; everything
# else (styles, untouched cell text, empty cells) is left exactly as-is.
# ---------------------------------------------------------------------------

# ============================== Overview ==================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws1.Range("A3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-23 00:45:57"

# Hyperlinks: display text tracks the new A/B column text, underlying
# relationship target stays pinned to the same cell slot (rId2 -> B2, rId3 -> B3).
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md", "", "", "e2e\fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md", "", "", "e2e\36038c19-52ce-4a1b-8036-de19daaeacb8.md")

# ================================ zh-cn ====================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws2.Range("G2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.zh-cn.xlf"
$ws2.Range("I2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws2.Range("J2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.zh-cn.xlf"

$ws2.Range("A3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("G3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-23 00:45:52"
$ws2.Range("I3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws2.Range("J3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.zh-cn.xlf"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7e3a148e73b7a69fc10850ed9049df5f22e7903/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md."

$ws2.Columns.Item(16).ColumnWidth = 39.17

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md", "", "", "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7eb3b2f47e49ddbaa6de1baf3d49f7a06430282/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md", "", "", "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md", "", "", "36038c19-52ce-4a1b-8036-de19daaeacb8.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7eb3b2f47e49ddbaa6de1baf3d49f7a06430282/e2e/fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md", "", "", "36038c19-52ce-4a1b-8036-de19daaeacb8.md")

# ================================ de-de ====================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws3.Range("G2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.de-de.xlf"
$ws3.Range("I2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md"
$ws3.Range("J2").Value = "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.3180788ccb1d5c858ef0ef8e59d53f0fa210ab48.de-de.xlf"

$ws3.Range("A3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("G3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-23 00:45:57"
$ws3.Range("I3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.md"
$ws3.Range("J3").Value = "36038c19-52ce-4a1b-8036-de19daaeacb8.4ad38bf24ece5dede3f5c797292c36779823b837.de-de.xlf"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7e3a148e73b7a69fc10850ed9049df5f22e7903/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md."

$ws3.Columns.Item(16).ColumnWidth = 39.17

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md", "", "", "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a1314b6a475308526d99cfd7a745a81f41bd4d96/e2e/36038c19-52ce-4a1b-8036-de19daaeacb8.md", "", "", "fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c664e18cfcefa78c8e3f5de7989a7b424ee7775b/e2e/fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md", "", "", "36038c19-52ce-4a1b-8036-de19daaeacb8.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a1314b6a475308526d99cfd7a745a81f41bd4d96/e2e/fbebf677-6ca6-45f1-aaa0-7f2528069a3b.md", "", "", "36038c19-52ce-4a1b-8036-de19daaeacb8.md")
